# GarysPosture App User Manual - v1.1 -> v1.2 update
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title paragraph: " v1.1" -> " v1.2", add <w:ind w:left="0"/>, and drop
#    the _GoBack bookmark that currently sits there (it gets moved to the
#    new "V1.1 - Initial Release" line we add at the end of the document).
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.First
$titlePara.Range.ParagraphFormat.LeftIndent = 0

$d.Content.Find.Execute("v1.1", $true, $false, $false, $false, $false, $true, 1, $false, "v1.2", 2) | Out-Null

try {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
} catch {
}

# ---------------------------------------------------------------------------
# 2. "Each session is saved..." paragraph gets reworded.
# ---------------------------------------------------------------------------
$oldSaved = "Each session is saved to the phones SDCard whenever you tap the Stop button or if you forget and exit the app. If you watch the phone when you tap Stop it will display a message for a few seconds telling you where the file is stored."
$newSaved = "Each session is saved to the phones SDCard as an Excel friendly csv formatted file.  The timestamp on each data point is the phone clock time in hours:minutes:seconds."
$d.Content.Find.Execute($oldSaved, $true, $false, $false, $false, $false, $true, 1, $false, $newSaved, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. "Add an email export of the data." gets reworded.
# ---------------------------------------------------------------------------
$oldEmail = "Add an email export of the data."
$newEmail = "Add an email/dropbox/Google Drive option to export the data."
$d.Content.Find.Execute($oldEmail, $true, $false, $false, $false, $false, $true, 1, $false, $newEmail, 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Append the new Changelog section at the end of the document.
# ---------------------------------------------------------------------------

# -- two blank paragraphs -----------------------------------------------
$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.RemoveNumbers()
$p.Style = "Normal"

$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.RemoveNumbers()
$p.Style = "Normal"

# -- page break paragraph -------------------------------------------------
$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.RemoveNumbers()
$p.Style = "Normal"
$p.Range.ParagraphFormat.LineSpacingRule = 0
$p.Range.ParagraphFormat.LineSpacing = 12
$pr = $p.Range
$pr.Collapse(0)
$pr.InsertAfter([char]12)

# -- Changelog heading ------------------------------------------------------
$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.RemoveNumbers()
$p.Style = "Heading1"
$pr = $p.Range
$pr.Collapse(0)
$pr.InsertAfter("Changelog")

# -- V1.2 changelog line -----------------------------------------------
$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.RemoveNumbers()
$p.Style = "Normal"
$pr = $p.Range
$pr.Collapse(0)
$pr.InsertAfter("V1.2 " + [char]0x2013 + " ")
$pr = $d.Content
$pr.Collapse(0)
$pr.InsertAfter([char]9)
$pr = $d.Content
$pr.Collapse(0)
$pr.InsertAfter("Timetamp on csv data is actual time")

# -- Data saved during monitoring line -----------------------------------
$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.RemoveNumbers()
$p.Style = "Normal"
$pr = $p.Range
$pr.Collapse(0)
$pr.InsertAfter([char]9)
$pr = $d.Content
$pr.Collapse(0)
$pr.InsertAfter("Data is saved to file during monitoring rather than at the end of a session")

# -- V1.1 Initial Release line, with the relocated _GoBack bookmark -------
$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.RemoveNumbers()
$p.Style = "Normal"
$pr = $p.Range
$pr.Collapse(0)
$pr.InsertAfter("V1.1 " + [char]0x2013 + " Initial ReleaseTEMPBOOKMARKPAD")

$bmRange = $d.Content
$bmRange.Find.Execute("V1.1 " + [char]0x2013 + " Initial Release", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)

$padRange = $d.Content
$padRange.Find.Execute("TEMPBOOKMARKPAD", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$padRange.Text = ""

Write-Output "edits applied"
